$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(65, 1).Value = 'Can you export image files?'
$ws.Cells.Item(65, 2).Value = 'llama3.2:latest'
$ws.Cells.Item(65, 3).Value = 'Yes, you can export image files from GEO. To do this:
1. Go to the Exporttab.
2. Select HTML Using... and choose your preferred image format (e.g., EMF).
3. An ExportTodialog box will open where you can enter a Filename and select the directory where it will be stored.
Note that exporting to HTML outputs an image of your log using spliced web-supported images, which are outputted to a location specified by the user at the point of saving, keeping them together.'

$ws.Cells.Item(66, 1).Value = 'Can you import your own data? Which format of data can be imported?'
$ws.Cells.Item(66, 2).Value = 'llama3.2:latest'
$ws.Cells.Item(66, 3).Value = 'Yes, you can import your own data into GEO. The preferred formats for importing data are comma-separated (CSV) or tab-delimited text files. These formats are more structured and preferred by vendors over space-delimited data. You can load these types of files without additional manipulation via the Data Import Wizard.'

$ws.Cells.Item(67, 1).Value = 'How to plot a curve in GEO?'
$ws.Cells.Item(67, 2).Value = 'llama3.2:latest'
$ws.Cells.Item(67, 3).Value = 'To plot a curve in GEO, you can use the Curve tool. This tool allows you to create and edit curves, which can be used to represent various types of data such as pressure, temperature, or other parameters.
You can also import existing curve data from other sources, such as CSV files or other GEO documents.
Additionally, you can use the Curve Editor to modify the properties of your curve, such as its color, line style, and marker type.
Please refer to the help section accessed through GEOGraph for more information on how to create and annotate a graph.'

$ws.Cells.Item(68, 1).Value = 'How to add title to the files -> setting up custom text?'
$ws.Cells.Item(68, 2).Value = 'llama3.2:latest'
$ws.Cells.Item(68, 3).Value = 'To add a title to a file in GEO, you can use the Format Text context tab. Highlight the text you wish to edit and select the format you want to implement from the available options. You can then draw the text smaller and lower than regular text by selecting the "Draw text smaller" option.
Additionally, if you need to insert superscript or subscript characters in your track text entry or free format text entry, highlight the text you wish to change and use either the Format Text context tab or the Rich Edit control.'

$ws.Cells.Item(69, 1).Value = 'How to open GEO?'
$ws.Cells.Item(69, 2).Value = 'llama3.2:latest'
$ws.Cells.Item(69, 3).Value = 'To open GEO, follow these steps:
1. Click on the Windows button and select "Geologix GEO Suite 8" from the list.
2. Then, click on "GEO".
3. In Windows Explorer, locate the installation location of GEO and double-click the "Geo.exe" file.
Alternatively, you can also open GEO by clicking on its icon in the Windows menu (if it appears) or by searching for it in the Start menu (if it''s not visible).'

# Widen column C to match target width (561 in OOXML units; ColumnWidth omits the ~0.8333 padding Excel adds)
$ws.Columns.Item(3).ColumnWidth = 560.1666666666667

# Restore default (non-custom) row heights on the newly added rows
$ws.Range("A65:A69").EntireRow.AutoFit()
